# Regenerate s_val columns (TB, d2S, K, IP, sum) for fairbanks_pete 2024 sheet
# after filtering save games out of the underlying computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B=0.01514828764759746; C=0.002777888934908601; D=3.900430680208489; E=0.496779210170732; G=4.415136066961727},
    @{Row=3; B=0.01514828764759746; C=0.04240448674262143; D=0.8054896365839992; E=0.496779210170732; G=1.35982162114495},
    @{Row=4; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=5; B=1.459612070389937; C=1.667794583268128; D=26.21740644021617; E=0.496779210170732; G=29.84159230404497},
    @{Row=6; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=7; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=8; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=9; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=10; B=1.459612070389937; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=4.429675500412797},
    @{Row=11; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=12; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=13; B=0.6753301551942219; C=0.3127903958511391; D=0.8054896365839992; E=0.496779210170732; G=2.290389397800092},
    @{Row=14; B=0.6753301551942219; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=2.997429241610044},
    @{Row=15; B=1.459612070389937; C=0.3127903958511391; D=0.1575252929769615; E=0.496779210170732; G=2.42670696938877},
    @{Row=16; B=1.459612070389937; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=4.429675500412797},
    @{Row=17; B=0.127881588408715; C=0.3127903958511391; D=0.8054896365839992; E=0.496779210170732; G=1.742940831014585},
    @{Row=18; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=19; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=20; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=21; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=22; B=0.6753301551942219; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=3.645393585217082},
    @{Row=23; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=24; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=25; B=1.459612070389937; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=4.429675500412797},
    @{Row=26; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=27; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=28; B=1.459612070389937; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=3.781711156805759},
    @{Row=29; B=0.6753301551942219; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=2.997429241610044},
    @{Row=30; B=0.6753301551942219; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=2.997429241610044},
    @{Row=31; B=0.6753301551942219; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=3.645393585217082},
    @{Row=32; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=33; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=34; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=35; B=3.230985683306322; C=1.667794583268128; D=0.8054896365839992; E=0.496779210170732; G=6.201049113329182},
    @{Row=36; B=0.3048080303191223; C=0.04240448674262143; D=0.1575252929769615; E=0.496779210170732; G=1.001517020209437},
    @{Row=37; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=38; B=1.459612070389937; C=1.667794583268128; D=26.21740644021617; E=0.496779210170732; G=29.84159230404497},
    @{Row=39; B=0.6753301551942219; C=0.3127903958511391; D=3.900430680208489; E=0.496779210170732; G=5.385330441424582},
    @{Row=40; B=3.230985683306322; C=1.667794583268128; D=0.1575252929769615; E=0.496779210170732; G=5.553084769722144},
    @{Row=41; B=0.6753301551942219; C=0.3127903958511391; D=0.8054896365839992; E=0.496779210170732; G=2.290389397800092},
    @{Row=42; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671},
    @{Row=43; B=3.230985683306322; C=1.667794583268128; D=3.900430680208489; E=0.496779210170732; G=9.295990156953671}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("G$r").Value = $item.G
}

